# Generate Report for handback
#
# Mirrors the localization "handback" report refresh:
#   - The shared status text "Ready for handoff" becomes
#     "Handed back: in sync with en-US" everywhere it is used (Overview
#     sheet B/C columns, and the Status column on the per-locale sheets).
#   - Each per-locale sheet (zh-cn, de-de) gains two new columns of data for
#     the already-handed-off rows: "Latest Target File" (E) and
#     "Latest Handback File" (F), populated with hyperlinks mirroring the
#     existing "Source File Name" / "Latest Handoff File" links.
#   - The "Latest Handback DateTime" column (G) is stamped with the new
#     handback timestamp for the rows that now have a handback.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$mdName = "41f7d767-a907-42ee-b2c9-5f80e94a7af9.md"

# ---------------------------------------------------------------------
# Overview sheet: the "Ready for handoff" status shown for the two real
# files (rows 2 and 3) flips to the new handed-back status in both the
# zh-cn and de-de columns.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

function Update-LocaleSheet($SheetName, $XlfName, $HandbackDateTime) {
    $ws = $wb.Worksheets.Item($SheetName)

    $mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/df63f1282318ce13f6f80adcb6b0ebbdfe5ee978/e2e/$mdName"
    $xlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/25339a55a6af8c4b2827a2fa626829806962babb/ol-handoff/OpenLocalizationTestOrg/oltest.$SheetName/xinjiang/ht/$XlfName"

    # Row 2 ------------------------------------------------------------
    $ws.Range("B2").Value = $newStatus

    $ws.Hyperlinks.Add($ws.Range("E2"), $mdUrl, "", $mdName, $mdName) | Out-Null
    $ws.Range("E2").Style = "HyperLink"

    $ws.Hyperlinks.Add($ws.Range("F2"), $xlfUrl, "", $XlfName, $XlfName) | Out-Null
    $ws.Range("F2").Style = "HyperLink"

    $ws.Range("G2").Value = $HandbackDateTime

    # Row 3 ------------------------------------------------------------
    $ws.Range("B3").Value = $newStatus

    $ws.Hyperlinks.Add($ws.Range("E3"), $mdUrl, "", $mdName, $mdName) | Out-Null
    $ws.Range("E3").Style = "HyperLink"

    $ws.Hyperlinks.Add($ws.Range("F3"), $xlfUrl, "", $XlfName, $XlfName) | Out-Null
    $ws.Range("F3").Style = "HyperLink"

    $ws.Range("G3").Value = $HandbackDateTime
}

Update-LocaleSheet "zh-cn" "41f7d767-a907-42ee-b2c9-5f80e94a7af9.c97fe14d0c784df19f1b84be13b2da20b3a6025c.zh-cn.xlf" "2016-02-16 10:23:45"
Update-LocaleSheet "de-de" "41f7d767-a907-42ee-b2c9-5f80e94a7af9.c97fe14d0c784df19f1b84be13b2da20b3a6025c.de-de.xlf" "2016-02-16 10:24:12"

Write-Output "Handback report generated"
